$wb = $excel.ActiveWorkbook

$wsSemilla9 = $wb.Worksheets.Item("Semilla 9")
$wsSemilla8 = $wb.Worksheets.Item("Semilla 8")

# Update data values (new use cases) on "Semilla 8" sheet, rows 9-14, columns C (MSIDN) and D (MSI)
$wsSemilla8.Range("C9").Value = "3045981670"
$wsSemilla8.Range("D9").Value = "732111324707276"

$wsSemilla8.Range("C10").Value = "3043209863"
$wsSemilla8.Range("D10").Value = "732111324707277"

$wsSemilla8.Range("C11").Value = "3045981670"
$wsSemilla8.Range("D11").Value = "732111324707276"

$wsSemilla8.Range("D12").Value = "732111193280551"
$wsSemilla8.Range("D13").Value = "732111193280544"
$wsSemilla8.Range("D14").Value = "732111193280535"

$wsSemilla8.Range("C12").Value = "3046010569"
$wsSemilla8.Range("C13").Value = "3046010523"
$wsSemilla8.Range("C14").Value = "3046008593"

# Update the selected cell/range on each sheet to match the saved view state
$wsSemilla9.Range("F2:G2").Select()
$wsSemilla8.Range("E16").Select()

$wsSemilla8.Activate()
